# Trading journal update: fill in Exit prices for several trades that have
# now closed, compute their realised P&L (column L), mark their Status
# (column T) as "Closed", and tag one of them with a Result (column K).
#
# Rows affected (all "Long" trades): 3, 5, 16, 17, 27, 35, 39

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3: CITY UNION BANK ---
$ws.Range("G3").Value = 197.5
$ws.Range("L3").Formula = "=(G3-D3)*H3"
$ws.Range("T3").Value = "Closed"

# --- Row 5: INDIAN BANK ---
$ws.Range("G5").Value = 669
$ws.Range("L5").Formula = "=(G5-D5)*H5"
$ws.Range("T5").Value = "Closed"

# --- Row 16: NAM-INDIA ---
$ws.Range("G16").Value = 796.8
$ws.Range("L16").Formula = "=(G16-D16)*H16"
$ws.Range("T16").Value = "Closed"
# This row previously used the red "pending" font (style index 3); once
# marked Closed the Status cell reverts to the normal default font, same
# as the rest of the already-closed rows (style index 1). Copy that look
# from a cell that already carries it.
$ws.Range("T3").Copy()
$ws.Range("T16").PasteSpecial(-4122)  # xlPasteFormats

# --- Row 17: ASTERDM ---
$ws.Range("G17").Value = 606
$ws.Range("L17").Formula = "=(G17-D17)*H17"
$ws.Range("T17").Value = "Closed"
$ws.Range("T3").Copy()
$ws.Range("T17").PasteSpecial(-4122)  # xlPasteFormats

$excel.CutCopyMode = 0

# --- Row 27: Supriya Lifescience ---
$ws.Range("G27").Value = 655.1
$ws.Range("K27").Value = "Loss"
$ws.Range("L27").Formula = "=(G27-D27)*H27"
$ws.Range("T27").Value = "Closed"

# --- Row 35: Mastek ---
$ws.Range("G35").Value = 2440
$ws.Range("L35").Formula = "=(G35-D35)*H35"
$ws.Range("T35").Value = "Closed"

# --- Row 39: Supriya Lifescience ---
$ws.Range("G39").Value = 648.44000000000005
$ws.Range("L39").Formula = "=(G39-D39)*H39"
$ws.Range("T39").Value = "Closed"

# Leave the cursor where the editor ended up.
$ws.Range("T35").Select() | Out-Null
